# Filter out low-score "Procedures" detail lines from the Patient Details
# document: remove the "flashlight", "sleep recording", and "." entries,
# while leaving "sleep" (and every other section/entry) untouched.

$d = $word.ActiveDocument

# Texts (first w:t run of each ListNumber paragraph) that should be dropped
# entirely, paragraph mark and all.
$removeTexts = @("flashlight", "sleep recording", ".")

# Walk the paragraph collection from the end towards the start so that
# deleting a paragraph never invalidates the index of one we still need to
# visit.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $styleName = $para.Style.NameLocal
    if ($styleName -ne "List Number") {
        continue
    }

    $fullText = $para.Range.Text
    # Each qualifying paragraph is built as: <entry><tab><tab><tab>" Page: N"
    # Split off everything from the first tab so we only compare the entry
    # label itself.
    $label = $fullText.Split("`t")[0].Trim()

    if ($removeTexts -contains $label) {
        $para.Range.Delete()
    }
}
